# Separate "dividends" and "tax" events: rename the Tax Withholding line
# items from "Withheld Tax on Dividends (...)" to "Tax Withholding (...)",
# correct the Foreign Currencies USD amounts that fed those events, and
# let the "Tax Withholding" sheet's column B width follow the (now
# shorter) label text.

$wb = $excel.ActiveWorkbook

# --- Tax Withholding sheet: shorten the per-symbol line-item labels ---
$wsTax = $wb.Worksheets.Item("Tax Withholding")
$wsTax.Range("B2").Value = "Tax Withholding (NVDA)"
$wsTax.Range("B3").Value = "Tax Withholding (APPL)"

# Column B no longer needs to be as wide for the shorter labels
# (was sized to fit "Withheld Tax on Dividends (...)", now just
# needs to fit "Tax Withholding (...)").
$wsTax.Columns.Item(2).ColumnWidth = 21.67

# --- Foreign Currencies sheet: corrected USD amounts behind the above ---
$wsFx = $wb.Worksheets.Item("Foreign Currencies")
$wsFx.Range("B2").Value = 1217.91
$wsFx.Range("B3").Value = 100
$wsFx.Range("B4").Value = 100
